$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct the missing underscore in the file name / path (row 7)
$ws.Range("A7").Value = "Overview_advis_codes_HL7_FHIR.pdf"
$ws.Range("E7").Value = "https://github.com/hl7dk/dk-medcom/blob/master/input/images/hospitalnotification/pdf/Overview_advis_codes_HL7_FHIR.pdf"

# Update the view state to match: no frozen/scrolled topLeftCell, selection moved to F15
$ws.Range("F15").Select()
